# Edit script: update open/close/high/low price, shares_outstanding, and fixed_ticker
# for all data rows (2-44) on the active worksheet, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
$rowData = @(
    @(2, 8.178509714471456, 6.507441520690918, 8.178509714471456, 6.507441520690918, 1319450062),
    @(3, 8.178509714471456, 6.507441520690918, 8.178509714471456, 6.507441520690918, 1319450062),
    @(4, 8.178509714471456, 6.507441520690918, 8.178509714471456, 6.507441520690918, 1319450062),
    @(5, 6.626802122327536, 6.569331645965576, 6.706376683858715, 5.800109917310798, 1319450062),
    @(6, 6.065172181419899, 5.887698173522949, 6.19827768734261, 5.330873619350438, 1319450062),
    @(7, 7.418213009081863, 8.224152565002441, 8.237511280641264, 6.848266226557223, 1319450062),
    @(8, 9.441327317043291, 9.593174934387209, 9.914733387647384, 9.213556623847344, 1319450062),
    @(9, 10.03732817027938, 10.65514850616455, 10.69991805313586, 9.634402982138988, 1319450062),
    @(10, 10.20317294788373, 10.24357223510742, 11.16827685228984, 9.875485666592835, 1319450062),
    @(11, 11.23116841498338, 11.30327892303467, 11.5136003763665, 11.098966063123, 1319450062),
    @(12, 10.58662531288694, 10.89427089691162, 10.94252870785988, 10.29707622258334, 1319450062),
    @(13, 10.89220583258957, 10.90001964569092, 11.18131095600954, 10.01317197817633, 1319450062),
    @(14, 12.88061590840274, 14.6006498336792, 15.02476690545092, 11.51008701761563, 1319450062),
    @(15, 13.39669352386079, 12.01681041717529, 14.05115230858778, 11.92218995152548, 1319450062),
    @(16, 12.24511180317164, 13.1265697479248, 13.96832149290186, 12.16570091825608, 1319450062),
    @(17, 12.18608641824067, 11.99417972564697, 13.20878960599863, 11.51841038301021, 1319450062),
    @(18, 12.64397718144804, 13.20002174377441, 13.67548054975801, 12.54324440962203, 1319450062),
    @(19, 12.94556268784594, 11.13561916351318, 13.04295867904445, 11.03010664620616, 1319450062),
    @(20, 11.72361531699597, 11.30637550354004, 12.0336813902634, 10.24282413935102, 1319450062),
    @(21, 13.27200134007031, 13.0494270324707, 14.500279507178, 12.98347926418032, 1319450062),
    @(22, 11.67852571750947, 10.62363719940186, 12.60882067318102, 10.01728455549422, 1319450062),
    @(23, 8.237117413551413, 8.15315055847168, 8.774502723606414, 7.237915201325159, 1319450062),
    @(24, 8.374201482234596, 8.204511642456055, 8.993570408857499, 7.737863976206329, 1319450062),
    @(25, 7.552380179380323, 9.485583305358888, 10.03547248950951, 7.427796404015487, 1319450062),
    @(26, 10.78702676211506, 12.63548755645752, 12.92186872468894, 10.57007128147795, 1319450062),
    @(27, 14.03357761567682, 13.96358394622803, 14.64601450533492, 13.66611397999722, 1319450062),
    @(28, 12.8631459008091, 13.6301736831665, 13.70246807020096, 12.58102097142668, 1319450062),
    @(29, 13.09375506209222, 12.75596714019775, 14.17823470693932, 12.54262628460936, 1319450062),
    @(30, 14.65556882223628, 14.26140880584717, 15.90971246496284, 13.84933289106825, 1319450062),
    @(31, 13.90197363220025, 14.07794857025146, 15.16086663180709, 12.21352287665498, 1319450062),
    @(32, 12.83826533064914, 12.37423229217529, 13.8300237650832, 12.31054080240857, 1319450062),
    @(33, 13.22088122899629, 15.4060001373291, 15.47944940938918, 12.45884437468464, 1319450062),
    @(34, 14.90940271365118, 14.43770217895508, 15.47359104146693, 14.29434279191638, 1319450062),
    @(35, 13.4251386817791, 13.45312690734863, 14.67995465661404, 12.73942048235315, 1319450062),
    @(36, 16.3606302533029, 15.96594524383545, 17.04193199245788, 15.22356260244737, 1319450062),
    @(37, 14.57417585273594, 16.00320243835449, 16.13569437176707, 14.52212583712397, 1319450062),
    @(38, 14.63156812893374, 14.52664756774902, 15.02263335974149, 14.07358456775807, 1319450062),
    @(39, 16.19901638802346, 16.95804405212402, 18.20707958476832, 15.62253859084705, 1319450062),
    @(40, 19.09029229993445, 18.72298622131348, 19.31260918551129, 16.05517806347853, 1319450062),
    @(41, 19.037042145779, 20.65793991088867, 22.02085924795754, 18.99323399770424, 1319450062),
    @(42, 20.19463907795179, 19.40605354309082, 21.52690790049732, 18.84767671848155, 1319450062),
    @(43, 16.19530774452335, 17.0643253326416, 17.95309186502507, 15.97311564054183, 1319450062),
    @(44, 20.18334157146687, 22.45135498046875, 23.18746476210288, 19.53178405586348, 1319450062)
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 4).Value = $entry[1]   # D - open_price
    $ws.Cells.Item($r, 5).Value = $entry[2]   # E - close_price
    $ws.Cells.Item($r, 6).Value = $entry[3]   # F - high_price
    $ws.Cells.Item($r, 7).Value = $entry[4]   # G - low_price
    $ws.Cells.Item($r, 8).Value = $entry[5]   # H - shares_outstanding
    $ws.Cells.Item($r, 9).Value = "HPE"       # I - fixed_ticker
}
